# Applies the "Language" rename + active-sheet/selection change to the
# Import001.xlsx workbook (Meta / Resources sheets).

$wb = $excel.ActiveWorkbook
$resourcesSheet = $wb.Worksheets.Item("Resources")

# Rename the "TwoLetterISOLanguageName" column header to "Language" on the
# Resources sheet (B1).
$resourcesSheet.Range("B1").Value = "Language"

# Resources sheet becomes the active sheet/tab, with B1 selected
# (Meta loses tabSelected automatically once Resources is activated).
$resourcesSheet.Activate()
$resourcesSheet.Range("B1").Select()
